$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "banana" requisito row (row 9), mirroring the formatting of
# the other rows: column A is italic (style used by most data rows),
# columns B/C use the plain bordered/wrap style (like row 4).
$ws.Range("A9").Value = "banana"
$ws.Range("B9").Value = "Funcional"
$ws.Range("C9").Value = "Tecnologia/Arquitetura"

$ws.Range("A9").Font.Italic = $true
$ws.Range("A9:C9").Borders.LineStyle = 1
$ws.Range("A9:C9").WrapText = $true

$ws.Range("D9").Select()
